# Append a new daily row (2025/09/30) to the tracking sheet, mirroring
# the existing rows' layout: date text, weekday text, hour number, rank number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A to be treated as text so the date-like string "2025/09/30"
# is stored literally instead of being auto-converted into a date serial
# number, then restore the default "Normal" style so the new row doesn't
# end up with a custom/visible number format (matching the other rows,
# which use the default style).
$ws.Range("A37").NumberFormat = "@"
$ws.Range("A37").Value = "2025/09/30"
$ws.Range("A37").Style = "Normal"

$ws.Range("B37").Value = "火"
$ws.Range("C37").Value = 1
$ws.Range("D37").Value = 15
